$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The first 16 data rows (rows 2 through 17) are removed; all subsequent
# rows shift up so that the data previously in row 18 becomes row 2, etc.
$ws.Range("A2:B17").EntireRow.Delete()
